$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 1249
$ws.Range("I94").Value = 1249
$ws.Range("K94").Value = 1249
$ws.Range("M94").Value = -798
$ws.Range("H106").Value = 4799.25
$ws.Range("J106").Value = 3496
$ws.Range("L106").Value = 3496
$ws.Range("N106").Value = -4758
$ws.Range("H112").Value = 9617344
$ws.Range("J112").Value = 9617344
$ws.Range("L112").Value = 28852032
$ws.Range("N112").Value = -28854248
$ws.Range("H137").Value = 6406
$ws.Range("I137").Value = 3612.9119
$ws.Range("J137").Value = 12737
$ws.Range("K137").Value = 10838.7357
$ws.Range("L137").Value = 38211
$ws.Range("M137").Value = -8288.735700000001
$ws.Range("N137").Value = -43311
$ws.Range("H138").Value = 6670240
$ws.Range("I138").Value = 728.6667
$ws.Range("J138").Value = 23820412
$ws.Range("K138").Value = 2186.0001
$ws.Range("L138").Value = 71461236
$ws.Range("M138").Value = 2953.9999
$ws.Range("N138").Value = -71471516

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2625.75
$ws.Range("I32").Value = 1305.3116
$ws.Range("J32").Value = 7046.3477
$ws.Range("K32").Value = 1305.3116
$ws.Range("L32").Value = 7046.3477
$ws.Range("M32").Value = -1018.3116
$ws.Range("N32").Value = -7620.3477
$ws.Range("H45").Value = 3547.923
$ws.Range("I45").Value = 3278.875
$ws.Range("J45").Value = 3978.4
$ws.Range("K45").Value = 3278.875
$ws.Range("L45").Value = 3978.4
$ws.Range("M45").Value = -2901.875
$ws.Range("N45").Value = -4732.4
$ws.Range("H74").Value = 62571176
$ws.Range("I74").Value = 83427576
$ws.Range("K74").Value = 83427576
$ws.Range("M74").Value = -83426702
$ws.Range("H77").Value = 62571176
$ws.Range("I77").Value = 83427576
$ws.Range("K77").Value = 417137880
$ws.Range("M77").Value = -417133512
$ws.Range("H82").Value = 44999.5
$ws.Range("J82").Value = 44999.5
$ws.Range("L82").Value = 44999.5
$ws.Range("N82").Value = -45721.5
$ws.Range("H85").Value = 44999.5
$ws.Range("J85").Value = 44999.5
$ws.Range("L85").Value = 44999.5
$ws.Range("N85").Value = -47495.5
$ws.Range("H97").Value = 1309.4722
$ws.Range("I97").Value = 1400
$ws.Range("K97").Value = 1400
$ws.Range("M97").Value = -904
$ws.Range("H110").Value = 15751.75
$ws.Range("I110").Value = 19433.21
$ws.Range("J110").Value = 1762.2
$ws.Range("K110").Value = 19433.21
$ws.Range("L110").Value = 1762.2
$ws.Range("M110").Value = -17388.21
$ws.Range("N110").Value = -5852.2
$ws.Range("H132").Value = 23313658
$ws.Range("I132").Value = 2080.75
$ws.Range("K132").Value = 6242.25
$ws.Range("M132").Value = -3712.25

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2659.7778
$ws.Range("I107").Value = 1857.3636
$ws.Range("J107").Value = 3920.7144
$ws.Range("K107").Value = 1857.3636
$ws.Range("L107").Value = 3920.7144
$ws.Range("M107").Value = 62.63640000000009
$ws.Range("N107").Value = -7760.7144
$ws.Range("H134").Value = 3032346.5
$ws.Range("I134").Value = 3127076.2
$ws.Range("K134").Value = 9381228.600000001
$ws.Range("M134").Value = -9378693.600000001
$ws.Range("H141").Value = 59499.75
$ws.Range("J141").Value = 53999.5
$ws.Range("L141").Value = 53999.5
$ws.Range("N141").Value = -64359.5

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 165
$ws.Range("J13").Value = 165
$ws.Range("L13").Value = 165
$ws.Range("N13").Value = -443
$ws.Range("H16").Value = 1034.4615
$ws.Range("I16").Value = 745.5
$ws.Range("J16").Value = 1496.8
$ws.Range("K16").Value = 745.5
$ws.Range("L16").Value = 1496.8
$ws.Range("M16").Value = -458.5
$ws.Range("N16").Value = -2070.8
$ws.Range("H26").Value = 16623.6
$ws.Range("I26").Value = 20404.5
$ws.Range("J26").Value = 1500
$ws.Range("K26").Value = 20404.5
$ws.Range("L26").Value = 1500
$ws.Range("M26").Value = -20117.5
$ws.Range("N26").Value = -2074
$ws.Range("H31").Value = 83339810
$ws.Range("I31").Value = 4113
$ws.Range("J31").Value = 200009780
$ws.Range("K31").Value = 4113
$ws.Range("L31").Value = 200009780
$ws.Range("M31").Value = -3818
$ws.Range("N31").Value = -200010370
$ws.Range("H34").Value = 83339810
$ws.Range("I34").Value = 4113
$ws.Range("J34").Value = 200009780
$ws.Range("K34").Value = 4113
$ws.Range("L34").Value = 200009780
$ws.Range("M34").Value = -3911
$ws.Range("N34").Value = -200010184
$ws.Range("H99").Value = 8870.957
$ws.Range("I99").Value = 5877.5713
$ws.Range("K99").Value = 5877.5713
$ws.Range("M99").Value = -4379.5713
$ws.Range("H113").Value = 1034.4615
$ws.Range("I113").Value = 745.5
$ws.Range("J113").Value = 1496.8
$ws.Range("K113").Value = 745.5
$ws.Range("L113").Value = 1496.8
$ws.Range("M113").Value = 1424.5
$ws.Range("N113").Value = -5836.8
$ws.Range("H126").Value = 8870.957
$ws.Range("I126").Value = 5877.5713
$ws.Range("K126").Value = 17632.7139
$ws.Range("M126").Value = -15162.7139
$ws.Range("H132").Value = 53582.7
$ws.Range("I132").Value = 65159.5
$ws.Range("K132").Value = 195478.5
$ws.Range("M132").Value = -192948.5
$ws.Range("H134").Value = 3783.7273
$ws.Range("I134").Value = 3912.1
$ws.Range("K134").Value = 11736.3
$ws.Range("M134").Value = -9201.299999999999
$ws.Range("H141").Value = 126369.125
$ws.Range("I141").Value = 35000
$ws.Range("J141").Value = 139421.86
$ws.Range("K141").Value = 35000
$ws.Range("L141").Value = 139421.86
$ws.Range("M141").Value = -29820
$ws.Range("N141").Value = -149781.86

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 724.75
$ws.Range("J92").Value = 999.5
$ws.Range("L92").Value = 2998.5
$ws.Range("N92").Value = -5494.5
$ws.Range("H128").Value = 99266.336
$ws.Range("I128").Value = 99266.336
$ws.Range("K128").Value = 297799.008
$ws.Range("M128").Value = -292819.008
$ws.Range("H133").Value = 7343.25
$ws.Range("J133").Value = 7000
$ws.Range("L133").Value = 21000
$ws.Range("N133").Value = -31120
$ws.Range("H134").Value = 1643.6666
$ws.Range("I134").Value = 1643.6666
$ws.Range("K134").Value = 4930.9998
$ws.Range("M134").Value = 139.0002000000004
$ws.Range("H137").Value = 2938.2727
$ws.Range("J137").Value = 3406.4
$ws.Range("L137").Value = 10219.2
$ws.Range("N137").Value = -20419.2

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2433.0278
$ws.Range("I102").Value = 1969.0344
$ws.Range("K102").Value = 1969.0344
$ws.Range("M102").Value = -347.0344
$ws.Range("H132").Value = 3174.2693
$ws.Range("I132").Value = 2518.1904
$ws.Range("K132").Value = 7554.5712
$ws.Range("M132").Value = -5024.5712

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3875.158
$ws.Range("I7").Value = 3248
$ws.Range("J7").Value = 4737.5
$ws.Range("K7").Value = 3248
$ws.Range("L7").Value = 4737.5
$ws.Range("M7").Value = -3136
$ws.Range("N7").Value = -4961.5
$ws.Range("H16").Value = 1733.375
$ws.Range("J16").Value = 1566.3334
$ws.Range("L16").Value = 1566.3334
$ws.Range("N16").Value = -1906.3334
$ws.Range("H22").Value = 1734.7812
$ws.Range("I22").Value = 1083.1765
$ws.Range("K22").Value = 1083.1765
$ws.Range("M22").Value = -788.1765
$ws.Range("H27").Value = 1734.7812
$ws.Range("I27").Value = 1083.1765
$ws.Range("K27").Value = 1083.1765
$ws.Range("M27").Value = -976.1765
$ws.Range("H93").Value = 1966.5454
$ws.Range("I93").Value = 1429.5454
$ws.Range("K93").Value = 1429.5454
$ws.Range("M93").Value = -181.5454
$ws.Range("H100").Value = 3565.5527
$ws.Range("I100").Value = 2996.0715
$ws.Range("J100").Value = 5160.1
$ws.Range("K100").Value = 2996.0715
$ws.Range("L100").Value = 5160.1
$ws.Range("M100").Value = -2455.0715
$ws.Range("N100").Value = -6242.1
$ws.Range("H122").Value = 3657.121
$ws.Range("I122").Value = 3076.5293
$ws.Range("J122").Value = 4274
$ws.Range("K122").Value = 9229.5879
$ws.Range("L122").Value = 12822
$ws.Range("M122").Value = -6779.5879
$ws.Range("N122").Value = -17722
$ws.Range("H126").Value = 3875.158
$ws.Range("I126").Value = 3248
$ws.Range("J126").Value = 4737.5
$ws.Range("K126").Value = 9744
$ws.Range("L126").Value = 14212.5
$ws.Range("M126").Value = -7274
$ws.Range("N126").Value = -19152.5
$ws.Range("I132").Value = 4061.5264
$ws.Range("K132").Value = 12184.5792
$ws.Range("M132").Value = -9654.5792

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2779.6553
$ws.Range("I81").Value = 2404.9167
$ws.Range("J81").Value = 4578.4
$ws.Range("K81").Value = 4809.8334
$ws.Range("L81").Value = 9156.799999999999
$ws.Range("M81").Value = -3748.8334
$ws.Range("N81").Value = -11278.8
$ws.Range("H84").Value = 2779.6553
$ws.Range("I84").Value = 2404.9167
$ws.Range("J84").Value = 4578.4
$ws.Range("K84").Value = 24049.167
$ws.Range("L84").Value = 45784
$ws.Range("M84").Value = -18745.167
$ws.Range("N84").Value = -56392
$ws.Range("H113").Value = 636.375
$ws.Range("I113").Value = 551.1539
$ws.Range("J113").Value = 737.0909
$ws.Range("K113").Value = 1653.4617
$ws.Range("L113").Value = 2211.2727
$ws.Range("M113").Value = 516.5382999999999
$ws.Range("N113").Value = -6551.2727
$ws.Range("H132").Value = 3768
$ws.Range("I132").Value = 3516.3333
$ws.Range("J132").Value = 4221
$ws.Range("K132").Value = 10548.9999
$ws.Range("L132").Value = 12663
$ws.Range("M132").Value = -8018.999899999999
$ws.Range("N132").Value = -17723
$ws.Range("H136").Value = 10063
$ws.Range("I136").Value = 6249.5
$ws.Range("J136").Value = 11334.167
$ws.Range("K136").Value = 18748.5
$ws.Range("L136").Value = 34002.501
$ws.Range("M136").Value = -16198.5
$ws.Range("N136").Value = -39102.501
